# Forecast update: last observation 2025-01-21.
# Replace the 19 (location) x 11 (weekly forecast column) block B1:L19
# with the refreshed VCI3M forecast values (row 1 = new forecast-week
# start dates, rows 2-19 = per-location forecast series).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: forecast week-start dates (Excel serial numbers)
$dates = @(45678,45685,45692,45699,45706,45713,45720,45727,45734,45741,45748)

# Rows 2-19: forecast values per location, one row per array
$forecastRows = @(
    @(43.2,42.5,42.5,43.2,44.4,46.3,48.6,51.5,54.8,58.4,62.2),
    @(48.7,49.3,50.5,52.2,54.4,57,60,63.3,67,70.9,74.8),
    @(34.9,34.5,34.5,34.9,35.7,36.8,38.2,39.9,42,44.4,47),
    @(36.5,34.8,33.6,33.2,33.5,34.5,36.4,39,42.4,46.4,50.7),
    @(51.3,51.2,51.4,52,53,54.3,56,58.1,60.7,63.8,67.1),
    @(25.5,24.1,23.8,24.6,26.4,28.8,31.8,35.2,38.7,42.2,45.6),
    @(43.4,42.7,42.6,43.1,44.2,45.8,48.1,50.9,54.3,58.1,62.2),
    @(44.9,44.7,45.3,46.5,48.2,50.2,52.4,54.7,57.2,59.7,62.1),
    @(32.8,30.9,29.6,28.8,28.4,28.6,29.2,30.4,32.1,34.3,37.1),
    @(54.9,54.8,55.1,55.8,57,58.8,61.3,64.5,68.2,72.3,76.6),
    @(44.9,44.7,45,45.8,46.9,48.4,50.2,52.4,54.8,57.4,60.1),
    @(57.2,57.6,58.2,58.9,59.6,60.3,61.3,62.5,64.1,66.1,68.6),
    @(47.5,46.1,45.4,45.4,46.2,47.9,50.4,53.7,57.5,61.8,66.2),
    @(50.3,50.7,51.2,51.8,52.3,52.8,53.1,53.4,53.7,54.1,54.6),
    @(35.9,34.1,33.4,33.6,34.7,36.7,39.4,42.9,46.9,51.4,56),
    @(39.5,37.2,35.6,34.9,35.1,36.3,38.5,41.7,45.5,50,54.8),
    @(45.7,45.4,46.5,49,52.8,57.6,63.4,69.8,76.8,84,91.1),
    @(39.8,39.5,40.7,43.2,46.8,51.5,56.9,62.9,69.2,75.7,82)
)

$data = New-Object 'object[,]' 19,11
for ($c = 0; $c -lt $dates.Length; $c++) {
    $data[0, $c] = $dates[$c]
}
for ($r = 0; $r -lt $forecastRows.Length; $r++) {
    $row = $forecastRows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $data[$r + 1, $c] = $row[$c]
    }
}

$ws.Range("B1:L19").Value = $data
